$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Sheet is protected; unprotect, make the edit, then restore protection.
$ws.Unprotect()

# Insert a new row above row 13 so the existing "(4) No more than 2000
# records may be entered." line shifts down to row 14, and write the new
# instruction line into the freshly inserted row 13.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "(4) Dates must be on or after January 2nd, 2018."

$ws.Protect($null, $true, $true, $true)
